$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Round row 5's numeric values down to 2 decimal places (custom accuracy).
$row5 = @{
    "B5" = 16.81;  "C5" = 12.3;   "D5" = 1.11;   "E5" = 36.54;  "F5" = 29.65;
    "G5" = 13.23;  "H5" = 49.96;  "I5" = 20.36;  "J5" = 8.98;   "K5" = 13.23;
    "L5" = 14.66;  "M5" = 15.44;  "N5" = 4.23;   "O5" = 13.16;  "P5" = 18.67;
    "Q5" = 11.18;  "R5" = 0.81;   "S5" = 0.74;   "T5" = 192.96; "U5" = 36.78;
    "V5" = 12.15;  "W5" = 24.62;  "X5" = 12.89;  "Y5" = 2.06;   "Z5" = 24.51;
    "AA5" = 10.73; "AB5" = 9.54;  "AC5" = 11.24; "AD5" = 15.39; "AE5" = 0.55;
    "AF5" = 45.47; "AG5" = 6.79;  "AH5" = 15.19
}
foreach ($addr in $row5.Keys) {
    $ws.Range($addr).Value = $row5[$addr]
}

# 2) Drop the last data row (row 6) entirely - "데이터 1000개" trim.
$ws.Rows.Item(6).Delete()

# 3) Narrow column K (11) from 8 to 7 characters wide.
$ws.Columns.Item(11).ColumnWidth = 6.17
